$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 228.83333
$ws.Range("I58").Value = 228.83333
$ws.Range("K58").Value = 686.49999
$ws.Range("M58").Value = -536.49999
$ws.Range("H88").Value = 1759.6154
$ws.Range("J88").Value = 1762.875
$ws.Range("L88").Value = 1762.875
$ws.Range("N88").Value = -2574.875
$ws.Range("H91").Value = 1759.6154
$ws.Range("J91").Value = 1762.875
$ws.Range("L91").Value = 1762.875
$ws.Range("N91").Value = -4570.875
$ws.Range("H112").Value = 5697.0454
$ws.Range("I112").Value = 18137.5
$ws.Range("J112").Value = 2932.5
$ws.Range("K112").Value = 54412.5
$ws.Range("L112").Value = 8797.5
$ws.Range("M112").Value = -53304.5
$ws.Range("N112").Value = -11013.5
$ws.Range("H137").Value = 5936.291
$ws.Range("I137").Value = 2212.7954
$ws.Range("K137").Value = 6638.3862
$ws.Range("M137").Value = -4088.3862

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7356819.5
$ws.Range("I32").Value = 11366755
$ws.Range("J32").Value = 5270.0835
$ws.Range("K32").Value = 11366755
$ws.Range("L32").Value = 5270.0835
$ws.Range("M32").Value = -11366468
$ws.Range("N32").Value = -5844.0835
$ws.Range("H45").Value = 4472
$ws.Range("I45").Value = 4444
$ws.Range("K45").Value = 4444
$ws.Range("M45").Value = -4067
$ws.Range("H61").Value = 723898
$ws.Range("I61").Value = 2914.1538
$ws.Range("K61").Value = 2914.1538
$ws.Range("M61").Value = -2702.1538
$ws.Range("H74").Value = 22967.178
$ws.Range("I74").Value = 1951.75
$ws.Range("J74").Value = 50987.75
$ws.Range("K74").Value = 1951.75
$ws.Range("L74").Value = 50987.75
$ws.Range("M74").Value = -1077.75
$ws.Range("N74").Value = -52735.75
$ws.Range("H77").Value = 22967.178
$ws.Range("I77").Value = 1951.75
$ws.Range("J77").Value = 50987.75
$ws.Range("K77").Value = 9758.75
$ws.Range("L77").Value = 254938.75
$ws.Range("M77").Value = -5390.75
$ws.Range("N77").Value = -263674.75
$ws.Range("H110").Value = 8059.2666
$ws.Range("I110").Value = 11168.9
$ws.Range("K110").Value = 11168.9
$ws.Range("M110").Value = -9123.9
$ws.Range("H132").Value = 5304082
$ws.Range("I132").Value = 2404.182
$ws.Range("K132").Value = 7212.545999999999
$ws.Range("M132").Value = -4682.545999999999
$ws.Range("H136").Value = 723898
$ws.Range("I136").Value = 2914.1538
$ws.Range("K136").Value = 8742.4614
$ws.Range("M136").Value = -6192.4614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 38945.668
$ws.Range("J81").Value = 38945.668
$ws.Range("L81").Value = 38945.668
$ws.Range("N81").Value = -41067.668
$ws.Range("H84").Value = 38945.668
$ws.Range("J84").Value = 38945.668
$ws.Range("L84").Value = 116837.004
$ws.Range("N84").Value = -127445.004
$ws.Range("H99").Value = 11050.714
$ws.Range("J99").Value = 5666.3335
$ws.Range("L99").Value = 5666.3335
$ws.Range("N99").Value = -8662.333500000001
$ws.Range("H105").Value = 2099.6
$ws.Range("I105").Value = 2099.6
$ws.Range("K105").Value = 2099.6
$ws.Range("M105").Value = -352.5999999999999
$ws.Range("H107").Value = 1384
$ws.Range("I107").Value = 1585
$ws.Range("K107").Value = 1585
$ws.Range("M107").Value = 335
$ws.Range("H110").Value = 34000
$ws.Range("J110").Value = 34000
$ws.Range("L110").Value = 34000
$ws.Range("N110").Value = -42180
$ws.Range("H134").Value = 46973.742
$ws.Range("I134").Value = 59720.53
$ws.Range("J134").Value = 25304.2
$ws.Range("K134").Value = 179161.59
$ws.Range("L134").Value = 75912.60000000001
$ws.Range("M134").Value = -176626.59
$ws.Range("N134").Value = -80982.60000000001
$ws.Range("H135").Value = 49000
$ws.Range("J135").Value = 49000
$ws.Range("L135").Value = 49000
$ws.Range("N135").Value = -59140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 15596.2
$ws.Range("J58").Value = 27696.385
$ws.Range("L58").Value = 27696.385
$ws.Range("N58").Value = -28102.385
$ws.Range("H132").Value = 42895010
$ws.Range("I132").Value = 2650.5264
$ws.Range("K132").Value = 7951.5792
$ws.Range("M132").Value = -5421.5792
$ws.Range("I134").Value = 2697.258
$ws.Range("J134").Value = 58833880
$ws.Range("K134").Value = 8091.773999999999
$ws.Range("L134").Value = 176501640
$ws.Range("M134").Value = -5556.773999999999
$ws.Range("N134").Value = -176506710
$ws.Range("H136").Value = 15596.2
$ws.Range("J136").Value = 27696.385
$ws.Range("L136").Value = 83089.155
$ws.Range("N136").Value = -88189.155

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 2045
$ws.Range("I115").Value = 1726.6666
$ws.Range("J115").Value = 3000
$ws.Range("K115").Value = 5179.9998
$ws.Range("L115").Value = 9000
$ws.Range("M115").Value = -4004.9998
$ws.Range("N115").Value = -11350
$ws.Range("H137").Value = 7871.8
$ws.Range("J137").Value = 14915.5
$ws.Range("L137").Value = 44746.5
$ws.Range("N137").Value = -54946.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6220.5386
$ws.Range("I102").Value = 6322.25
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 6322.25
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -4700.25
$ws.Range("N102").Value = -8244
$ws.Range("H122").Value = 2509.0715
$ws.Range("I122").Value = 2534
$ws.Range("J122").Value = 2484.1428
$ws.Range("K122").Value = 7602
$ws.Range("L122").Value = 7452.428400000001
$ws.Range("M122").Value = -5152
$ws.Range("N122").Value = -12352.4284
$ws.Range("H123").Value = 55236.875
$ws.Range("J123").Value = 55236.875
$ws.Range("L123").Value = 55236.875
$ws.Range("N123").Value = -60136.875
$ws.Range("H132").Value = 856012.2
$ws.Range("I132").Value = 4319.9287
$ws.Range("J132").Value = 2843294.2
$ws.Range("K132").Value = 12959.7861
$ws.Range("L132").Value = 8529882.600000001
$ws.Range("M132").Value = -10429.7861
$ws.Range("N132").Value = -8534942.600000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2951.7778
$ws.Range("I82").Value = 2816.1428
$ws.Range("K82").Value = 2816.1428
$ws.Range("M82").Value = -2455.1428
$ws.Range("H85").Value = 2951.7778
$ws.Range("I85").Value = 2816.1428
$ws.Range("K85").Value = 2816.1428
$ws.Range("M85").Value = -1568.1428
$ws.Range("H123").Value = 21583.334
$ws.Range("J123").Value = 21583.334
$ws.Range("L123").Value = 21583.334
$ws.Range("N123").Value = -31383.334
$ws.Range("H132").Value = 2504307.8
$ws.Range("I132").Value = 4021.6667
$ws.Range("J132").Value = 3688653.8
$ws.Range("K132").Value = 12065.0001
$ws.Range("L132").Value = 11065961.4
$ws.Range("M132").Value = -9535.000100000001
$ws.Range("N132").Value = -11071021.4
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H135").Value = 30342.9
$ws.Range("J135").Value = 30342.9
$ws.Range("L135").Value = 30342.9
$ws.Range("N135").Value = -40482.9
$ws.Range("H136").Value = 1070779.1
$ws.Range("I136").Value = 21421
$ws.Range("J136").Value = 1678302.2
$ws.Range("K136").Value = 64263
$ws.Range("L136").Value = 5034906.6
$ws.Range("M136").Value = -61713
$ws.Range("N136").Value = -5040006.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 10000
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H62").Value = 18626.5
$ws.Range("I62").Value = 28013.4
$ws.Range("J62").Value = 9239.6
$ws.Range("K62").Value = 28013.4
$ws.Range("L62").Value = 9239.6
$ws.Range("M62").Value = -27389.4
$ws.Range("N62").Value = -10487.6
$ws.Range("H65").Value = 18626.5
$ws.Range("I65").Value = 28013.4
$ws.Range("J65").Value = 9239.6
$ws.Range("K65").Value = 140067
$ws.Range("L65").Value = 46198
$ws.Range("M65").Value = -136947
$ws.Range("N65").Value = -52438
$ws.Range("H74").Value = 42686.25
$ws.Range("J74").Value = 42686.25
$ws.Range("L74").Value = 42686.25
$ws.Range("N74").Value = -44558.25
$ws.Range("H77").Value = 42686.25
$ws.Range("J77").Value = 42686.25
$ws.Range("L77").Value = 128058.75
$ws.Range("N77").Value = -137418.75
$ws.Range("H122").Value = 3935.0667
$ws.Range("I122").Value = 2259.0557
$ws.Range("J122").Value = 6449.0835
$ws.Range("K122").Value = 6777.1671
$ws.Range("L122").Value = 19347.2505
$ws.Range("M122").Value = -4327.1671
$ws.Range("N122").Value = -24247.2505
